$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Sunday, Jan 08"
$ws.Range("C9").Value = "10:20 PM"
$ws.Range("D9").Value = "FR6628"
$ws.Range("E9").Value = "London"
$ws.Range("F9").Value = "(LTN)"
$ws.Range("G9").Value = "Ryanair "
$ws.Range("H9").Value = "B738"
$ws.Range("I9").Value = "(EI-EBZ)"
$ws.Range("J9").Value = "10:10 PM"
$ws.Range("L9").Value = "0 hours, -10 minutes"
$ws.Range("K9").Borders.LineStyle = -4142
$ws.Range("M9").Borders.LineStyle = -4142
